$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

# Row 28 <-> Row 29: coin swap (Kaspa <-> Monero) with updated figures
Set-TextValue 'B28' 'Monero'
Set-TextValue 'C28' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D28' '159.32'
Set-TextValue 'E28' '  -3.80%  '
Set-TextValue 'B29' 'Kaspa'
Set-TextValue 'C29' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D29' '0.135'
Set-TextValue 'E29' '  +19.06%  '

# Remaining price/volume updates
Set-TextValue 'D2' '36.182.38'
Set-TextValue 'E2' '  -3.78%  '
Set-TextValue 'D3' '1.977.44'
Set-TextValue 'E3' '  -3.19%  '
Set-TextValue 'E4' '  +0.14%  '
Set-TextValue 'D5' '243.84'
Set-TextValue 'E5' '  -3.93%  '
Set-TextValue 'D6' '0.626'
Set-TextValue 'E6' '  -3.25%  '
Set-TextValue 'D7' '62.16'
Set-TextValue 'E7' '  -1.59%  '
Set-TextValue 'E8' '  +0.03%  '
Set-TextValue 'D9' '0.374'
Set-TextValue 'E9' '  -0.28%  '
Set-TextValue 'D10' '56.71'
Set-TextValue 'E10' '  -3.83%  '
Set-TextValue 'D11' '0.0796'
Set-TextValue 'E11' '  +5.47%  '
Set-TextValue 'E12' '  -0.94%  '
Set-TextValue 'D13' '0.864'
Set-TextValue 'E13' '  -5.20%  '
Set-TextValue 'D14' '22.45'
Set-TextValue 'E14' '  +9.83%  '
Set-TextValue 'D15' '14.04'
Set-TextValue 'E15' '  -6.98%  '
Set-TextValue 'D16' '2.274.84'
Set-TextValue 'E16' '  -2.84%  '
Set-TextValue 'D17' '5.43'
Set-TextValue 'E17' '  -2.66%  '
Set-TextValue 'D18' '1.991.17'
Set-TextValue 'E18' '  -2.42%  '
Set-TextValue 'D19' '36.087.61'
Set-TextValue 'E19' '  -3.80%  '
Set-TextValue 'D20' '71.06'
Set-TextValue 'E20' '  -3.07%  '
Set-TextValue 'D21' '0.0₃0866'
Set-TextValue 'E21' '  -0.75%  '
Set-TextValue 'D22' '238.43'
Set-TextValue 'E22' '  +0.71%  '
Set-TextValue 'D23' '5.25'
Set-TextValue 'E23' '  -1.75%  '
Set-TextValue 'D24' '0.999'
Set-TextValue 'E24' '  -0.19%  '
Set-TextValue 'E25' '  -10.17%  '
Set-TextValue 'D26' '2.30'
Set-TextValue 'E26' '  -1.50%  '
Set-TextValue 'D27' '9.73'
Set-TextValue 'E27' '  +1.95%  '
Set-TextValue 'D30' '19.77'
Set-TextValue 'E30' '  -0.21%  '
Set-TextValue 'D31' '0.119'
Set-TextValue 'E31' '  -1.95%  '
Set-TextValue 'D32' '4.90'
Set-TextValue 'E32' '  -5.32%  '
Set-TextValue 'D33' '1.14'
Set-TextValue 'E33' '  -6.13%  '
Set-TextValue 'D34' '0.0620'
Set-TextValue 'E34' '  +0.84%  '
Set-TextValue 'D35' '4.38'
Set-TextValue 'E35' '  -6.44%  '
Set-TextValue 'D36' '6.32'
Set-TextValue 'E36' '  +4.80%  '
Set-TextValue 'E37' '  +0.31%  '
Set-TextValue 'E38' '  -6.98%  '
Set-TextValue 'D39' '1.84'
Set-TextValue 'E39' '  +1.99%  '
Set-TextValue 'D40' '3.12'
Set-TextValue 'E40' '  +14.79%  '
Set-TextValue 'D41' '0.0990'
Set-TextValue 'E41' '  -6.20%  '
Set-TextValue 'D42' '1.23'
Set-TextValue 'E42' '  -0.71%  '
Set-TextValue 'D43' '0.0213'
Set-TextValue 'E43' '  -2.80%  '
Set-TextValue 'D44' '2.84'
Set-TextValue 'E44' '  -3.25%  '
Set-TextValue 'D45' '1.09'
Set-TextValue 'E45' '  -4.01%  '
Set-TextValue 'D46' '92.99'
Set-TextValue 'E46' '  -2.43%  '
Set-TextValue 'D47' '16.17'
Set-TextValue 'E47' '  -4.32%  '
Set-TextValue 'D48' '7.53'
Set-TextValue 'E48' '  -7.00%  '
Set-TextValue 'D49' '1.350.48'
Set-TextValue 'E49' '  -5.53%  '
Set-TextValue 'D50' '2.85'
Set-TextValue 'E50' '  -3.31%  '
Set-TextValue 'D51' '2.168.73'
Set-TextValue 'E51' '  -2.60%  '
